$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Film title / Note(score) pairs for rows 2-35 (row 5 "Un p'tit truc en plus" is unchanged)
$data = @(
    @("Furiosa - Une saga Mad Max", "7.6"),
    @("Civil War", "7"),
    @("Le Deuxième Acte", "6.3"),
    @("Un p'tit truc en plus", "6.8"),
    @("The Fall Guy", "6.3"),
    @("Challengers", "6.6"),
    @("Atlas", "4.6"),
    @("La Planète des singes - Le Nouveau Royaume", "6.4"),
    @("Marcello Mio", "5.6"),
    @("Comme un lundi", "6.7"),
    @("When Evil Lurks", "6.5"),
    @("Les Trois Fantastiques", "6.4"),
    @("Borgo", "7"),
    @("Arthur the King", "6.5"),
    @("South Park - La fin de l'obésité", "7.1"),
    @("Jusqu’au bout du monde", "6.5"),
    @("The Beekeeper", "5.2"),
    @("L'Idée d'être avec toi", "5.2"),
    @("Blue & Compagnie", "6.1"),
    @("Border Line", "6.8"),
    @("Adagio", "7"),
    @("Moi aussi", "5.1"),
    @("Monkey Man", "6.5"),
    @("Le Tableau volé", "6.1"),
    @("Mon oni à moi", "5.9"),
    @("Back to Black", "5.7"),
    @("Chien blanc", "5.7"),
    @("Rebel Moon : Partie 2 - L'Entailleuse", "4"),
    @("Les Intrus", "4.2"),
    @("Heroico", "6.5"),
    @("Unfrosted - L'épopée de la Pop-Tart", "4.5"),
    @("Baghead", "4.9"),
    @("Baby Ruby", "5.4"),
    @("La Mère de la mariée", "4")
)

$startRow = 2
$endRow = $startRow + $data.Count - 1

# Column B ("Note") must stay text (as in the source data) instead of being auto-coerced to a
# number, so mark the range as text before writing, then restore the default "Normal" style so
# no stray number-format sticks around on the cells afterwards.
$noteRange = $ws.Range("B$startRow`:B$endRow")
$noteRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
}

$noteRange.Style = "Normal"

# The refreshed list is shorter than before (34 films instead of 36) - drop the now-stale trailing rows
$ws.Rows("36:37").Delete()
